$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00373831775700935
$ws.Range("C2").Value = 0.0205607476635514
$ws.Range("D2").Value = 0.923364485981308
$ws.Range("E2").Value = 0.0598130841121495
$ws.Range("F2").Value = 0.934579439252336
$ws.Range("G2").Value = 0.0560747663551402
$ws.Range("H2").Value = 0.0168224299065421
$ws.Range("I2").Value = 0.364485981308411
$ws.Range("J2").Value = 0.0186915887850467
$ws.Range("K2").Value = 0.0168224299065421
$ws.Range("L2").Value = 0.00373831775700935
$ws.Range("M2").Value = 0.011214953271028
$ws.Range("N2").Value = 0.998130841121495
$ws.Range("O2").Value = 0.00560747663551402
$ws.Range("P2").Value = 0.00934579439252336
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0.951401869158878
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0.011214953271028
$ws.Range("U2").Value = 0.00186915887850467
$ws.Range("V2").Value = 0.00934579439252336
$ws.Range("W2").Value = 0.0897196261682243
$ws.Range("X2").Value = 0.0411214953271028
$ws.Range("B3").Value = 0.0691588785046729
$ws.Range("C3").Value = 0.902803738317757
$ws.Range("D3").Value = 0.0504672897196262
$ws.Range("E3").Value = 0.0130841121495327
$ws.Range("F3").Value = 0.00934579439252336
$ws.Range("G3").Value = 0.91588785046729
$ws.Range("H3").Value = 0.906542056074766
$ws.Range("I3").Value = 0.624299065420561
$ws.Range("J3").Value = 0.0728971962616822
$ws.Range("K3").Value = 0.0242990654205607
$ws.Range("L3").Value = 0.00560747663551402
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0.0130841121495327
$ws.Range("Q3").Value = 0.94392523364486
$ws.Range("R3").Value = 0.0411214953271028
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0.0598130841121495
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0.00747663551401869
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0.00186915887850467
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0.00186915887850467
$ws.Range("D4").Value = 0.0168224299065421
$ws.Range("E4").Value = 0.914018691588785
$ws.Range("F4").Value = 0.0542056074766355
$ws.Range("G4").Value = 0.00934579439252336
$ws.Range("H4").Value = 0.00186915887850467
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0.00747663551401869
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0.177570093457944
$ws.Range("M4").Value = 0.0149532710280374
$ws.Range("N4").Value = 0.00186915887850467
$ws.Range("O4").Value = 0.994392523364486
$ws.Range("P4").Value = 0.00373831775700935
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0.00373831775700935
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = 0.011214953271028
$ws.Range("U4").Value = 0.998130841121495
$ws.Range("V4").Value = 0.00934579439252336
$ws.Range("W4").Value = 0.906542056074766
$ws.Range("X4").Value = 0.955140186915888
$ws.Range("B5").Value = 0.927102803738318
$ws.Range("C5").Value = 0.0728971962616822
$ws.Range("D5").Value = 0.00934579439252336
$ws.Range("E5").Value = 0.0130841121495327
$ws.Range("F5").Value = 0.00186915887850467
$ws.Range("G5").Value = 0.0186915887850467
$ws.Range("H5").Value = 0.0747663551401869
$ws.Range("I5").Value = 0.011214953271028
$ws.Range("J5").Value = 0.900934579439252
$ws.Range("K5").Value = 0.958878504672897
$ws.Range("L5").Value = 0.813084112149533
$ws.Range("M5").Value = 0.973831775700935
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0.973831775700935
$ws.Range("Q5").Value = 0.0560747663551402
$ws.Range("R5").Value = 0.00373831775700935
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0.917757009345794
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 0.973831775700935
$ws.Range("W5").Value = 0.00186915887850467
$ws.Range("X5").Value = 0.00186915887850467
